$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tmp = $wb.Worksheets.Add()
$tmp.Name = "TmpHelper"
$ws = $wb.Worksheets.Item("Sheet1")

# Helper cell carrying a Text ("@") number format, used (once) to force
# numeric-looking values to be stored as text/shared-strings like the source data.
$textHelper = $tmp.Cells.Item(1,1)
$textHelper.NumberFormat = "@"
$textHelper.Value2 = "x"

# Helper preserving the existing bold/bordered "column A" style (style index 1)
$aStyleHelper = $tmp.Cells.Item(1,2)
$ws.Cells.Item(1,1).Copy() | Out-Null
$aStyleHelper.PasteSpecial(-4122) | Out-Null

function Set-PlainValue($cell, $val) {
    $cell.Value2 = $val
}

function Set-TextValue($cell, $val) {
    $textHelper.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

# Wipe the old (4-row) table completely before laying out the new one.
$ws.Cells.Clear()

# Row 1
Set-PlainValue $ws.Cells.Item(1,1) 'Year'
Set-PlainValue $ws.Cells.Item(1,2) 'Total business travel (‘000km) 32,039'
Set-PlainValue $ws.Cells.Item(1,3) 'Individual business travel (‘000km per FTE) 6.9'
Set-PlainValue $ws.Cells.Item(1,4) 'Of which by train (%) 38'
Set-PlainValue $ws.Cells.Item(1,5) 'Of which by air (%) 16'
Set-PlainValue $ws.Cells.Item(1,6) 'Of which by private car (%) 2'
Set-PlainValue $ws.Cells.Item(1,7) 'Of which by company car (%) 40'
Set-PlainValue $ws.Cells.Item(1,8) 'Of which by rental car (%) 3'
Set-PlainValue $ws.Cells.Item(1,11) 'The carbon emissions from commuter travel were included for the first time in 2019, based on statistical figures.'
Set-PlainValue $ws.Cells.Item(1,12) 'They amounted to 3,249 tonnes and are based on data supplied by the German Federal Statistical Office for the'
Set-PlainValue $ws.Cells.Item(1,13) 'distribution of commuter traffic and the calculation methods of the VfU.'
Set-PlainValue $ws.Cells.Item(1,14) '3.3 GRI 303: Water and effluents 2018'
Set-PlainValue $ws.Cells.Item(1,16) '303-1 Interactions with water as a shared resource'
Set-PlainValue $ws.Cells.Item(1,17) 'DZ BANK AG only uses drinking water provided by regional suppliers. Our water consumption does not have'
Set-PlainValue $ws.Cells.Item(1,18) 'any significant impact on water sources.'
Set-PlainValue $ws.Cells.Item(1,20) '303-2 Management of water discharge related impacts'
Set-PlainValue $ws.Cells.Item(1,21) 'We generally only produce waste water that is comparable with household effluents. It is not reused or recycled.'
Set-PlainValue $ws.Cells.Item(1,23) '303-5 Water consumption'
Set-PlainValue $ws.Cells.Item(1,24) 'WATER CONSUMPTION'
Set-PlainValue $ws.Cells.Item(1,26) 'cubic meters (m3) 2019'
Set-PlainValue $ws.Cells.Item(1,27) 'Total water consumption 101,057'
Set-PlainValue $ws.Cells.Item(1,28) 'Individual water consumption (m3 per FTE) 22'

# Row 2
Set-TextValue $ws.Cells.Item(2,1) '2018'
Set-TextValue $ws.Cells.Item(2,2) '32,766'
Set-TextValue $ws.Cells.Item(2,3) '6.9'
Set-TextValue $ws.Cells.Item(2,4) '39'
Set-TextValue $ws.Cells.Item(2,5) '16'
Set-TextValue $ws.Cells.Item(2,6) '3'
Set-TextValue $ws.Cells.Item(2,7) '39'
Set-TextValue $ws.Cells.Item(2,8) '3'
Set-TextValue $ws.Cells.Item(2,26) '2018'
Set-TextValue $ws.Cells.Item(2,27) '96,293'
Set-TextValue $ws.Cells.Item(2,28) '20'

# Row 3
Set-TextValue $ws.Cells.Item(3,1) '2017'
Set-TextValue $ws.Cells.Item(3,2) '36,644'
Set-TextValue $ws.Cells.Item(3,3) '7.7'
Set-TextValue $ws.Cells.Item(3,4) '42'
Set-TextValue $ws.Cells.Item(3,5) '16'
Set-TextValue $ws.Cells.Item(3,6) '3'
Set-TextValue $ws.Cells.Item(3,7) '37'
Set-TextValue $ws.Cells.Item(3,8) '3'
Set-TextValue $ws.Cells.Item(3,26) '2017'
Set-TextValue $ws.Cells.Item(3,27) '94,022'
Set-TextValue $ws.Cells.Item(3,28) '20'

# Re-apply the bold/border "column A" style to A1:A3 (style index 1, reused - no new style created).
$aStyleHelper.Copy() | Out-Null
$ws.Range("A1:A3").PasteSpecial(-4122) | Out-Null

$wb.Worksheets.Item("TmpHelper").Delete() | Out-Null

$ws.Range("A1").Select() | Out-Null
